# Generate Report for Handback
# Updates the status of the 7bc79fcf-000e-47b3-a00c-ff4582d3354f.md file
# from "Ready for handoff" to "Handed back: in sync with en-US" across the
# Overview, zh-cn and de-de sheets, and records new "Latest Handback
# DateTime" values for the zh-cn and de-de handback reports.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusText
$wsZhCn.Range("G3").Value = "2016-02-17 06:13:23"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusText
$wsDeDe.Range("G3").Value = "2016-02-17 06:13:40"
